$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.188941
$ws.Range("H2").Value = 18.566823
$ws.Range("I2").Value = 0.5329955127131565
$ws.Range("J2").Value = 0.5329955127131565
$ws.Range("M2").Value = 6.188941
$ws.Range("N2").Value = 18.566823
$ws.Range("O2").Value = 0.5329955127131565
$ws.Range("P2").Value = 0.5329955127131565
$ws.Range("Q2").Value = 38.302990701481
$ws.Range("R2").Value = 344.726916313329
$ws.Range("S2").Value = 0.2840842165723605
$ws.Range("T2").Value = 0.2840842165723605
$ws.Range("G3").Value = 6.188941
$ws.Range("H3").Value = 18.566823
$ws.Range("I3").Value = 0.5329955127131565
$ws.Range("J3").Value = 0.5329955127131565
$ws.Range("N3").Value = 7.126554
$ws.Range("O3").Value = 0.2045811124018362
$ws.Range("P3").Value = 0.2045811124018361
$ws.Range("Q3").Value = 14.701940746438
$ws.Range("R3").Value = 132.317466717942
$ws.Range("S3").Value = 0.1090408148960446
$ws.Range("T3").Value = 0.1090408148960445
$ws.Range("G4").Value = 6.188941
$ws.Range("H4").Value = 18.566823
$ws.Range("I4").Value = 0.5329955127131565
$ws.Range("J4").Value = 0.5329955127131565
$ws.Range("M4").Value = 3.003166333333333
$ws.Range("N4").Value = 9.009499
$ws.Range("O4").Value = 0.2586345837838639
$ws.Range("P4").Value = 0.2586345837838639
$ws.Range("Q4").Value = 18.58641925018633
$ws.Range("R4").Value = 167.277773251677
$ws.Range("S4").Value = 0.1378510725892344
$ws.Range("T4").Value = 0.1378510725892344
$ws.Range("G5").Value = 6.188941
$ws.Range("H5").Value = 18.566823
$ws.Range("I5").Value = 0.5329955127131565
$ws.Range("J5").Value = 0.5329955127131565
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.043994
$ws.Range("N5").Value = 0.131982
$ws.Range("O5").Value = 0.003788791101143573
$ws.Range("P5").Value = 0.003788791101143573
$ws.Range("Q5").Value = 0.272276270354
$ws.Range("R5").Value = 2.450486433186
$ws.Range("S5").Value = 0.002019408655517064
$ws.Range("T5").Value = 0.002019408655517063
$ws.Range("H6").Value = 7.126554
$ws.Range("I6").Value = 0.2045811124018362
$ws.Range("J6").Value = 0.2045811124018361
$ws.Range("M6").Value = 6.188941
$ws.Range("N6").Value = 18.566823
$ws.Range("O6").Value = 0.5329955127131565
$ws.Range("P6").Value = 0.5329955127131565
$ws.Range("Q6").Value = 14.701940746438
$ws.Range("R6").Value = 132.317466717942
$ws.Range("S6").Value = 0.1090408148960446
$ws.Range("T6").Value = 0.1090408148960445
$ws.Range("H7").Value = 7.126554
$ws.Range("I7").Value = 0.2045811124018362
$ws.Range("J7").Value = 0.2045811124018361
$ws.Range("N7").Value = 7.126554
$ws.Range("O7").Value = 0.2045811124018362
$ws.Range("P7").Value = 0.2045811124018361
$ws.Range("Q7").Value = 5.643085768324
$ws.Range("R7").Value = 50.78777191491599
$ws.Range("S7").Value = 0.04185343155157272
$ws.Range("T7").Value = 0.04185343155157271
$ws.Range("H8").Value = 7.126554
$ws.Range("I8").Value = 0.2045811124018362
$ws.Range("J8").Value = 0.2045811124018361
$ws.Range("M8").Value = 3.003166333333333
$ws.Range("N8").Value = 9.009499
$ws.Range("O8").Value = 0.2586345837838639
$ws.Range("P8").Value = 0.2586345837838639
$ws.Range("Q8").Value = 7.134075681827333
$ws.Range("R8").Value = 64.20668113644599
$ws.Range("S8").Value = 0.05291175085608877
$ws.Range("T8").Value = 0.05291175085608876
$ws.Range("H9").Value = 7.126554
$ws.Range("I9").Value = 0.2045811124018362
$ws.Range("J9").Value = 0.2045811124018361
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.043994
$ws.Range("N9").Value = 0.131982
$ws.Range("O9").Value = 0.003788791101143573
$ws.Range("P9").Value = 0.003788791101143573
$ws.Range("Q9").Value = 0.104508538892
$ws.Range("R9").Value = 0.9405768500279998
$ws.Range("S9").Value = 0.0007751150981301299
$ws.Range("T9").Value = 0.0007751150981301297
$ws.Range("G10").Value = 3.003166333333333
$ws.Range("H10").Value = 9.009499
$ws.Range("I10").Value = 0.2586345837838639
$ws.Range("J10").Value = 0.2586345837838639
$ws.Range("M10").Value = 6.188941
$ws.Range("N10").Value = 18.566823
$ws.Range("O10").Value = 0.5329955127131565
$ws.Range("P10").Value = 0.5329955127131565
$ws.Range("Q10").Value = 18.58641925018633
$ws.Range("R10").Value = 167.277773251677
$ws.Range("S10").Value = 0.1378510725892344
$ws.Range("T10").Value = 0.1378510725892344
$ws.Range("G11").Value = 3.003166333333333
$ws.Range("H11").Value = 9.009499
$ws.Range("I11").Value = 0.2586345837838639
$ws.Range("J11").Value = 0.2586345837838639
$ws.Range("N11").Value = 7.126554
$ws.Range("O11").Value = 0.2045811124018362
$ws.Range("P11").Value = 0.2045811124018361
$ws.Range("Q11").Value = 7.134075681827333
$ws.Range("R11").Value = 64.20668113644599
$ws.Range("S11").Value = 0.05291175085608877
$ws.Range("T11").Value = 0.05291175085608876
$ws.Range("G12").Value = 3.003166333333333
$ws.Range("H12").Value = 9.009499
$ws.Range("I12").Value = 0.2586345837838639
$ws.Range("J12").Value = 0.2586345837838639
$ws.Range("M12").Value = 3.003166333333333
$ws.Range("N12").Value = 9.009499
$ws.Range("O12").Value = 0.2586345837838639
$ws.Range("P12").Value = 0.2586345837838639
$ws.Range("Q12").Value = 9.019008025666777
$ws.Range("R12").Value = 81.17107223100101
$ws.Range("S12").Value = 0.06689184792905251
$ws.Range("T12").Value = 0.06689184792905251
$ws.Range("G13").Value = 3.003166333333333
$ws.Range("H13").Value = 9.009499
$ws.Range("I13").Value = 0.2586345837838639
$ws.Range("J13").Value = 0.2586345837838639
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.043994
$ws.Range("N13").Value = 0.131982
$ws.Range("O13").Value = 0.003788791101143573
$ws.Range("P13").Value = 0.003788791101143573
$ws.Range("Q13").Value = 0.1321212996686666
$ws.Range("R13").Value = 1.189091697018
$ws.Range("S13").Value = 0.0009799124094882754
$ws.Range("T13").Value = 0.0009799124094882754
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.043994
$ws.Range("H14").Value = 0.131982
$ws.Range("I14").Value = 0.003788791101143573
$ws.Range("J14").Value = 0.003788791101143573
$ws.Range("M14").Value = 6.188941
$ws.Range("N14").Value = 18.566823
$ws.Range("O14").Value = 0.5329955127131565
$ws.Range("P14").Value = 0.5329955127131565
$ws.Range("Q14").Value = 0.272276270354
$ws.Range("R14").Value = 2.450486433186
$ws.Range("S14").Value = 0.002019408655517064
$ws.Range("T14").Value = 0.002019408655517063
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.043994
$ws.Range("H15").Value = 0.131982
$ws.Range("I15").Value = 0.003788791101143573
$ws.Range("J15").Value = 0.003788791101143573
$ws.Range("N15").Value = 7.126554
$ws.Range("O15").Value = 0.2045811124018362
$ws.Range("P15").Value = 0.2045811124018361
$ws.Range("Q15").Value = 0.104508538892
$ws.Range("R15").Value = 0.9405768500279998
$ws.Range("S15").Value = 0.0007751150981301299
$ws.Range("T15").Value = 0.0007751150981301297
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.043994
$ws.Range("H16").Value = 0.131982
$ws.Range("I16").Value = 0.003788791101143573
$ws.Range("J16").Value = 0.003788791101143573
$ws.Range("M16").Value = 3.003166333333333
$ws.Range("N16").Value = 9.009499
$ws.Range("O16").Value = 0.2586345837838639
$ws.Range("P16").Value = 0.2586345837838639
$ws.Range("Q16").Value = 0.1321212996686666
$ws.Range("R16").Value = 1.189091697018
$ws.Range("S16").Value = 0.0009799124094882754
$ws.Range("T16").Value = 0.0009799124094882754
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.043994
$ws.Range("H17").Value = 0.131982
$ws.Range("I17").Value = 0.003788791101143573
$ws.Range("J17").Value = 0.003788791101143573
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.043994
$ws.Range("N17").Value = 0.131982
$ws.Range("O17").Value = 0.003788791101143573
$ws.Range("P17").Value = 0.003788791101143573
$ws.Range("Q17").Value = 0.001935472036
$ws.Range("R17").Value = 0.017419248324
$ws.Range("S17").Value = 0.00001435493800810473
$ws.Range("T17").Value = 0.00001435493800810473

Write-Host "Applied TPM update to Cxadr-Cxadr sheet"
